# Scheduled-runner refresh of market price / profit columns (H-N)
# across several worksheets, matching the latest pulled price data.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(62, 8).Value = 3197.5  # H62: 2388.3333 -> 3197.5
$ws.Cells.Item(62, 9).Value = 1895  # I62: 2070.7144 -> 1895
$ws.Cells.Item(62, 10).Value = 4500  # J62: 3500 -> 4500
$ws.Cells.Item(62, 11).Value = 1895  # K62: 2070.7144 -> 1895
$ws.Cells.Item(62, 12).Value = 4500  # L62: 3500 -> 4500
$ws.Cells.Item(62, 13).Value = -1271  # M62: -1446.7144 -> -1271
$ws.Cells.Item(62, 14).Value = -5748  # N62: -4748 -> -5748
$ws.Cells.Item(64, 8).Value = 5999.8335  # H64: 0 -> 5999.8335
$ws.Cells.Item(64, 9).Value = 4999.8  # I64: 0 -> 4999.8
$ws.Cells.Item(64, 10).Value = 11000  # J64: 0 -> 11000
$ws.Cells.Item(64, 11).Value = 4999.8  # K64: 0 -> 4999.8
$ws.Cells.Item(64, 12).Value = 11000  # L64: 0 -> 11000
$ws.Cells.Item(64, 13).Value = -4751.8  # M64: (blank) -> -4751.8
$ws.Cells.Item(64, 14).Value = -11496  # N64: (blank) -> -11496
$ws.Cells.Item(65, 8).Value = 3197.5  # H65: 2388.3333 -> 3197.5
$ws.Cells.Item(65, 9).Value = 1895  # I65: 2070.7144 -> 1895
$ws.Cells.Item(65, 10).Value = 4500  # J65: 3500 -> 4500
$ws.Cells.Item(65, 11).Value = 9475  # K65: 10353.572 -> 9475
$ws.Cells.Item(65, 12).Value = 22500  # L65: 17500 -> 22500
$ws.Cells.Item(65, 13).Value = -6355  # M65: -7233.572 -> -6355
$ws.Cells.Item(65, 14).Value = -28740  # N65: -23740 -> -28740
$ws.Cells.Item(67, 8).Value = 5999.8335  # H67: 0 -> 5999.8335
$ws.Cells.Item(67, 9).Value = 4999.8  # I67: 0 -> 4999.8
$ws.Cells.Item(67, 10).Value = 11000  # J67: 0 -> 11000
$ws.Cells.Item(67, 11).Value = 4999.8  # K67: 0 -> 4999.8
$ws.Cells.Item(67, 12).Value = 11000  # L67: 0 -> 11000
$ws.Cells.Item(67, 13).Value = -4141.8  # M67: (blank) -> -4141.8
$ws.Cells.Item(67, 14).Value = -12716  # N67: (blank) -> -12716
$ws.Cells.Item(97, 8).Value = 3165  # H97: 3910 -> 3165
$ws.Cells.Item(97, 10).Value = 3165  # J97: 3910 -> 3165
$ws.Cells.Item(97, 12).Value = 9495  # L97: 11730 -> 9495
$ws.Cells.Item(97, 14).Value = -10487  # N97: -12722 -> -10487
$ws.Cells.Item(112, 8).Value = 5000  # H112: 5500 -> 5000
$ws.Cells.Item(112, 10).Value = 5000  # J112: 5600 -> 5000
$ws.Cells.Item(112, 12).Value = 15000  # L112: 16800 -> 15000
$ws.Cells.Item(112, 14).Value = -17216  # N112: -19016 -> -17216
$ws.Cells.Item(123, 8).Value = 85780  # H123: 0 -> 85780
$ws.Cells.Item(123, 10).Value = 85780  # J123: 0 -> 85780
$ws.Cells.Item(123, 12).Value = 85780  # L123: 0 -> 85780
$ws.Cells.Item(123, 14).Value = -95580  # N123: (blank) -> -95580
$ws.Cells.Item(138, 8).Value = 1739.6786  # H138: 2362.6191 -> 1739.6786
$ws.Cells.Item(138, 9).Value = 839.15  # I138: 1229.3846 -> 839.15
$ws.Cells.Item(138, 10).Value = 3991  # J138: 4204.125 -> 3991
$ws.Cells.Item(138, 11).Value = 2517.45  # K138: 3688.1538 -> 2517.45
$ws.Cells.Item(138, 12).Value = 11973  # L138: 12612.375 -> 11973
$ws.Cells.Item(138, 13).Value = 2622.55  # M138: 1451.8462 -> 2622.55
$ws.Cells.Item(138, 14).Value = -22253  # N138: -22892.375 -> -22253

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(61, 8).Value = 2909  # H61: 3947.25 -> 2909
$ws.Cells.Item(61, 9).Value = 3039  # I61: 4596.6665 -> 3039
$ws.Cells.Item(61, 11).Value = 3039  # K61: 4596.6665 -> 3039
$ws.Cells.Item(61, 13).Value = -2827  # M61: -4384.6665 -> -2827
$ws.Cells.Item(74, 8).Value = 2425.6924  # H74: 2137.375 -> 2425.6924
$ws.Cells.Item(74, 9).Value = 2425.6924  # I74: 2199.2 -> 2425.6924
$ws.Cells.Item(74, 10).Value = 0  # J74: 1210 -> 0
$ws.Cells.Item(74, 11).Value = 2425.6924  # K74: 2199.2 -> 2425.6924
$ws.Cells.Item(74, 12).Value = 0  # L74: 1210 -> 0
$ws.Cells.Item(74, 13).ClearContents()  # M74: -1325.2 -> (blank)
$ws.Cells.Item(74, 14).Value = -1551.6924  # N74: -2958 -> -1551.6924
$ws.Cells.Item(77, 8).Value = 2425.6924  # H77: 2137.375 -> 2425.6924
$ws.Cells.Item(77, 9).Value = 2425.6924  # I77: 2199.2 -> 2425.6924
$ws.Cells.Item(77, 10).Value = 0  # J77: 1210 -> 0
$ws.Cells.Item(77, 11).Value = 12128.462  # K77: 10996 -> 12128.462
$ws.Cells.Item(77, 12).Value = 0  # L77: 6050 -> 0
$ws.Cells.Item(77, 13).ClearContents()  # M77: -6628 -> (blank)
$ws.Cells.Item(77, 14).Value = -7760.462  # N77: -14786 -> -7760.462
$ws.Cells.Item(122, 8).Value = 1536.7142  # H122: 1154.7858 -> 1536.7142
$ws.Cells.Item(122, 9).Value = 1351.4  # I122: 1013.9167 -> 1351.4
$ws.Cells.Item(122, 11).Value = 4054.2  # K122: 3041.7501 -> 4054.2
$ws.Cells.Item(122, 13).Value = -1604.2  # M122: -591.7501000000002 -> -1604.2
$ws.Cells.Item(136, 8).Value = 2909  # H136: 3947.25 -> 2909
$ws.Cells.Item(136, 9).Value = 3039  # I136: 4596.6665 -> 3039
$ws.Cells.Item(136, 11).Value = 9117  # K136: 13789.9995 -> 9117
$ws.Cells.Item(136, 13).Value = -6567  # M136: -11239.9995 -> -6567

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(86, 8).Value = 23299.9  # H86: 38749.5 -> 23299.9
$ws.Cells.Item(86, 9).Value = 26398.8  # I86: 99999 -> 26398.8
$ws.Cells.Item(86, 10).Value = 20201  # J86: 18333 -> 20201
$ws.Cells.Item(86, 11).Value = 26398.8  # K86: 99999 -> 26398.8
$ws.Cells.Item(86, 12).Value = 20201  # L86: 18333 -> 20201
$ws.Cells.Item(86, 13).Value = -25275.8  # M86: -98876 -> -25275.8
$ws.Cells.Item(86, 14).Value = -22447  # N86: -20579 -> -22447
$ws.Cells.Item(89, 8).Value = 23299.9  # H89: 38749.5 -> 23299.9
$ws.Cells.Item(89, 9).Value = 26398.8  # I89: 99999 -> 26398.8
$ws.Cells.Item(89, 10).Value = 20201  # J89: 18333 -> 20201
$ws.Cells.Item(89, 11).Value = 131994  # K89: 499995 -> 131994
$ws.Cells.Item(89, 12).Value = 101005  # L89: 91665 -> 101005
$ws.Cells.Item(89, 13).Value = -126378  # M89: -494379 -> -126378
$ws.Cells.Item(89, 14).Value = -112237  # N89: -102897 -> -112237
$ws.Cells.Item(94, 8).Value = 851.3333  # H94: 1213.8948 -> 851.3333
$ws.Cells.Item(94, 9).Value = 857.6429000000001  # I94: 2021.5 -> 857.6429000000001
$ws.Cells.Item(94, 10).Value = 829.25  # J94: 998.5333000000001 -> 829.25
$ws.Cells.Item(94, 11).Value = 857.6429000000001  # K94: 2021.5 -> 857.6429000000001
$ws.Cells.Item(94, 12).Value = 829.25  # L94: 998.5333000000001 -> 829.25
$ws.Cells.Item(94, 13).Value = -406.6429000000001  # M94: -1570.5 -> -406.6429000000001
$ws.Cells.Item(94, 14).Value = -1731.25  # N94: -1900.5333 -> -1731.25
$ws.Cells.Item(99, 8).Value = 772.6  # H99: 843.7143 -> 772.6
$ws.Cells.Item(99, 9).Value = 772.6  # I99: 784.4 -> 772.6
$ws.Cells.Item(99, 10).Value = 0  # J99: 992 -> 0
$ws.Cells.Item(99, 11).Value = 772.6  # K99: 784.4 -> 772.6
$ws.Cells.Item(99, 12).Value = 0  # L99: 992 -> 0
$ws.Cells.Item(99, 13).ClearContents()  # M99: 713.6 -> (blank)
$ws.Cells.Item(99, 14).Value = 725.4  # N99: -3988 -> 725.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 2938.3333  # H31: 4548.5 -> 2938.3333
$ws.Cells.Item(31, 9).Value = 1650.6  # I31: 0 -> 1650.6
$ws.Cells.Item(31, 10).Value = 4548  # J31: 4548.5 -> 4548
$ws.Cells.Item(31, 11).Value = 1650.6  # K31: 0 -> 1650.6
$ws.Cells.Item(31, 12).Value = 4548  # L31: 4548.5 -> 4548
$ws.Cells.Item(31, 13).Value = -1355.6  # M31: (blank) -> -1355.6
$ws.Cells.Item(31, 14).Value = -5138  # N31: -5138.5 -> -5138
$ws.Cells.Item(34, 8).Value = 2938.3333  # H34: 4548.5 -> 2938.3333
$ws.Cells.Item(34, 9).Value = 1650.6  # I34: 0 -> 1650.6
$ws.Cells.Item(34, 10).Value = 4548  # J34: 4548.5 -> 4548
$ws.Cells.Item(34, 11).Value = 1650.6  # K34: 0 -> 1650.6
$ws.Cells.Item(34, 12).Value = 4548  # L34: 4548.5 -> 4548
$ws.Cells.Item(34, 13).Value = -1448.6  # M34: (blank) -> -1448.6
$ws.Cells.Item(34, 14).Value = -4952  # N34: -4952.5 -> -4952
$ws.Cells.Item(58, 8).Value = 2620.7144  # H58: 6166.1665 -> 2620.7144
$ws.Cells.Item(58, 9).Value = 2553.7273  # I58: 5249.25 -> 2553.7273
$ws.Cells.Item(58, 10).Value = 2866.3333  # J58: 8000 -> 2866.3333
$ws.Cells.Item(58, 11).Value = 2553.7273  # K58: 5249.25 -> 2553.7273
$ws.Cells.Item(58, 12).Value = 2866.3333  # L58: 8000 -> 2866.3333
$ws.Cells.Item(58, 13).Value = -2350.7273  # M58: -5046.25 -> -2350.7273
$ws.Cells.Item(58, 14).Value = -3272.3333  # N58: -8406 -> -3272.3333
$ws.Cells.Item(86, 8).Value = 3177.3635  # H86: 3381.6667 -> 3177.3635
$ws.Cells.Item(86, 9).Value = 3235  # I86: 3379.2 -> 3235
$ws.Cells.Item(86, 10).Value = 3023.6667  # J86: 3394 -> 3023.6667
$ws.Cells.Item(86, 11).Value = 3235  # K86: 3379.2 -> 3235
$ws.Cells.Item(86, 12).Value = 3023.6667  # L86: 3394 -> 3023.6667
$ws.Cells.Item(86, 13).Value = -2112  # M86: -2256.2 -> -2112
$ws.Cells.Item(86, 14).Value = -5269.6667  # N86: -5640 -> -5269.6667
$ws.Cells.Item(89, 8).Value = 3177.3635  # H89: 3381.6667 -> 3177.3635
$ws.Cells.Item(89, 9).Value = 3235  # I89: 3379.2 -> 3235
$ws.Cells.Item(89, 10).Value = 3023.6667  # J89: 3394 -> 3023.6667
$ws.Cells.Item(89, 11).Value = 16175  # K89: 16896 -> 16175
$ws.Cells.Item(89, 12).Value = 15118.3335  # L89: 16970 -> 15118.3335
$ws.Cells.Item(89, 13).Value = -10559  # M89: -11280 -> -10559
$ws.Cells.Item(89, 14).Value = -26350.3335  # N89: -28202 -> -26350.3335
$ws.Cells.Item(136, 8).Value = 2620.7144  # H136: 6166.1665 -> 2620.7144
$ws.Cells.Item(136, 9).Value = 2553.7273  # I136: 5249.25 -> 2553.7273
$ws.Cells.Item(136, 10).Value = 2866.3333  # J136: 8000 -> 2866.3333
$ws.Cells.Item(136, 11).Value = 7661.1819  # K136: 15747.75 -> 7661.1819
$ws.Cells.Item(136, 12).Value = 8598.999899999999  # L136: 24000 -> 8598.999899999999
$ws.Cells.Item(136, 13).Value = -5111.1819  # M136: -13197.75 -> -5111.1819
$ws.Cells.Item(136, 14).Value = -13698.9999  # N136: -29100 -> -13698.9999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(80, 8).Value = 2182.7144  # H80: 4998.3335 -> 2182.7144
$ws.Cells.Item(80, 9).Value = 2129.8333  # I80: 4998.3335 -> 2129.8333
$ws.Cells.Item(80, 10).Value = 2500  # J80: 0 -> 2500
$ws.Cells.Item(80, 11).Value = 2129.8333  # K80: 4998.3335 -> 2129.8333
$ws.Cells.Item(80, 12).Value = 2500  # L80: 0 -> 2500
$ws.Cells.Item(80, 13).Value = -1131.8333  # M80: -4000.3335 -> -1131.8333
$ws.Cells.Item(80, 14).Value = -4496  # N80: (blank) -> -4496
$ws.Cells.Item(83, 8).Value = 2182.7144  # H83: 4998.3335 -> 2182.7144
$ws.Cells.Item(83, 9).Value = 2129.8333  # I83: 4998.3335 -> 2129.8333
$ws.Cells.Item(83, 10).Value = 2500  # J83: 0 -> 2500
$ws.Cells.Item(83, 11).Value = 10649.1665  # K83: 24991.6675 -> 10649.1665
$ws.Cells.Item(83, 12).Value = 12500  # L83: 0 -> 12500
$ws.Cells.Item(83, 13).Value = -5657.166499999999  # M83: -19999.6675 -> -5657.166499999999
$ws.Cells.Item(83, 14).Value = -22484  # N83: (blank) -> -22484
$ws.Cells.Item(113, 8).Value = 4056.6667  # H113: 4048.5715 -> 4056.6667
$ws.Cells.Item(113, 9).Value = 4210  # I113: 4168 -> 4210
$ws.Cells.Item(113, 11).Value = 4210  # K113: 4168 -> 4210
$ws.Cells.Item(113, 13).Value = -2040  # M113: -1998 -> -2040
$ws.Cells.Item(122, 8).Value = 1836.875  # H122: 4612.6 -> 1836.875
$ws.Cells.Item(122, 9).Value = 1439.1  # I122: 4513.75 -> 1439.1
$ws.Cells.Item(122, 10).Value = 2499.8333  # J122: 5008 -> 2499.8333
$ws.Cells.Item(122, 11).Value = 4317.299999999999  # K122: 13541.25 -> 4317.299999999999
$ws.Cells.Item(122, 12).Value = 7499.499899999999  # L122: 15024 -> 7499.499899999999
$ws.Cells.Item(122, 13).Value = -1867.299999999999  # M122: -11091.25 -> -1867.299999999999
$ws.Cells.Item(122, 14).Value = -12399.4999  # N122: -19924 -> -12399.4999
$ws.Cells.Item(132, 8).Value = 2385.484  # H132: 2566.4285 -> 2385.484
$ws.Cells.Item(132, 9).Value = 1878.7391  # I132: 2044.5 -> 1878.7391
$ws.Cells.Item(132, 10).Value = 3842.375  # J132: 3871.25 -> 3842.375
$ws.Cells.Item(132, 11).Value = 5636.2173  # K132: 6133.5 -> 5636.2173
$ws.Cells.Item(132, 12).Value = 11527.125  # L132: 11613.75 -> 11527.125
$ws.Cells.Item(132, 13).Value = -3106.2173  # M132: -3603.5 -> -3106.2173
$ws.Cells.Item(132, 14).Value = -16587.125  # N132: -16673.75 -> -16587.125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(16, 8).Value = 1428.25  # H16: 1439 -> 1428.25
$ws.Cells.Item(16, 9).Value = 904.3333  # I16: 918.6667 -> 904.3333
$ws.Cells.Item(16, 11).Value = 904.3333  # K16: 918.6667 -> 904.3333
$ws.Cells.Item(16, 13).Value = -734.3333  # M16: -748.6667 -> -734.3333
$ws.Cells.Item(68, 8).Value = 2265.889  # H68: 4804.25 -> 2265.889
$ws.Cells.Item(68, 9).Value = 2342.2856  # I68: 5739 -> 2342.2856
$ws.Cells.Item(68, 10).Value = 1998.5  # J68: 2000 -> 1998.5
$ws.Cells.Item(68, 11).Value = 2342.2856  # K68: 5739 -> 2342.2856
$ws.Cells.Item(68, 12).Value = 1998.5  # L68: 2000 -> 1998.5
$ws.Cells.Item(68, 13).Value = -1593.2856  # M68: -4990 -> -1593.2856
$ws.Cells.Item(68, 14).Value = -3496.5  # N68: -3498 -> -3496.5
$ws.Cells.Item(71, 8).Value = 2265.889  # H71: 4804.25 -> 2265.889
$ws.Cells.Item(71, 9).Value = 2342.2856  # I71: 5739 -> 2342.2856
$ws.Cells.Item(71, 10).Value = 1998.5  # J71: 2000 -> 1998.5
$ws.Cells.Item(71, 11).Value = 11711.428  # K71: 28695 -> 11711.428
$ws.Cells.Item(71, 12).Value = 9992.5  # L71: 10000 -> 9992.5
$ws.Cells.Item(71, 13).Value = -7967.428  # M71: -24951 -> -7967.428
$ws.Cells.Item(71, 14).Value = -17480.5  # N71: -17488 -> -17480.5
$ws.Cells.Item(122, 8).Value = 2245.6155  # H122: 2266.3845 -> 2245.6155
$ws.Cells.Item(122, 9).Value = 1955.1666  # I122: 2000.1666 -> 1955.1666
$ws.Cells.Item(122, 11).Value = 5865.4998  # K122: 6000.4998 -> 5865.4998
$ws.Cells.Item(122, 13).Value = -3415.4998  # M122: -3550.4998 -> -3415.4998

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(122, 8).Value = 2785.4546  # H122: 2567.2 -> 2785.4546
$ws.Cells.Item(122, 9).Value = 2785.4546  # I122: 2296.3333 -> 2785.4546
$ws.Cells.Item(122, 10).Value = 0  # J122: 5005 -> 0
$ws.Cells.Item(122, 11).Value = 8356.363799999999  # K122: 6888.999899999999 -> 8356.363799999999
$ws.Cells.Item(122, 12).Value = 0  # L122: 15015 -> 0
$ws.Cells.Item(122, 13).ClearContents()  # M122: -4438.999899999999 -> (blank)
$ws.Cells.Item(122, 14).Value = -5906.363799999999  # N122: -19915 -> -5906.363799999999
$ws.Cells.Item(136, 8).Value = 7364.107  # H136: 8206.76 -> 7364.107
$ws.Cells.Item(136, 9).Value = 5703.8076  # I136: 6403.174 -> 5703.8076
$ws.Cells.Item(136, 11).Value = 17111.4228  # K136: 19209.522 -> 17111.4228
$ws.Cells.Item(136, 13).Value = -14561.4228  # M136: -16659.522 -> -14561.4228
